$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Recorded By" column (G) lists the person(s)/system that recorded a
# session, e.g. "System, dnasr281@gmail.com". This pass re-orders those
# two-part values so the human/email author is listed first and "System"
# is listed second (e.g. "dnasr281@gmail.com, System"), while leaving the
# "System, backup@backdoor.com" entries untouched.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notlike "System, *") { continue }

    $rest = $val.Substring(8)
    if ($rest -eq "" -or $rest.Contains(",")) { continue }
    if ($rest -eq "backup@backdoor.com") { continue }

    $cell.Value = "$rest, System"
}
